$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new worker (ADOL ANTONIO TORRES TRESPALACIOS, period 2111) is inserted as the
# first detail row (row 16). The existing worker (SILFREDO MARTINEZ CABARCAS) rows
# shift down by one and are re-sorted so that their periods run in ascending order
# (2207 .. 2404) instead of the previous descending order. Row/column layout, styles
# and merged cells stay the same; only the content of C:G for rows 16-38 changes.

$rows = @(
    @{ Row = 16; Doc = "1126118302"; Name = "ADOL ANTONIO TORRES TRESPALACIOS"; Period = "2111"; Mora = 30666;  Salario = 1000000 },
    @{ Row = 17; Doc = "73136692";   Name = "SILFREDO MARTINEZ CABARCAS";       Period = "2207"; Mora = 24000;  Salario = 1500000 },
    @{ Row = 18; Doc = "73136692";   Name = "SILFREDO MARTINEZ CABARCAS";       Period = "2208"; Mora = 60000;  Salario = 1500000 },
    @{ Row = 19; Doc = "73136692";   Name = "SILFREDO MARTINEZ CABARCAS";       Period = "2209"; Mora = 60000;  Salario = 1500000 },
    @{ Row = 20; Doc = "73136692";   Name = "SILFREDO MARTINEZ CABARCAS";       Period = "2210"; Mora = 60000;  Salario = 1500000 },
    @{ Row = 21; Doc = "73136692";   Name = "SILFREDO MARTINEZ CABARCAS";       Period = "2211"; Mora = 60000;  Salario = 1500000 },
    @{ Row = 22; Doc = "73136692";   Name = "SILFREDO MARTINEZ CABARCAS";       Period = "2212"; Mora = 60000;  Salario = 1500000 },
    @{ Row = 23; Doc = "73136692";   Name = "SILFREDO MARTINEZ CABARCAS";       Period = "2301"; Mora = 60000;  Salario = 1500000 },
    @{ Row = 24; Doc = "73136692";   Name = "SILFREDO MARTINEZ CABARCAS";       Period = "2302"; Mora = 60000;  Salario = 1500000 },
    @{ Row = 25; Doc = "73136692";   Name = "SILFREDO MARTINEZ CABARCAS";       Period = "2303"; Mora = 60000;  Salario = 1500000 },
    @{ Row = 26; Doc = "73136692";   Name = "SILFREDO MARTINEZ CABARCAS";       Period = "2304"; Mora = 60000;  Salario = 1500000 },
    @{ Row = 27; Doc = "73136692";   Name = "SILFREDO MARTINEZ CABARCAS";       Period = "2305"; Mora = 60000;  Salario = 1500000 },
    @{ Row = 28; Doc = "73136692";   Name = "SILFREDO MARTINEZ CABARCAS";       Period = "2306"; Mora = 60000;  Salario = 1500000 },
    @{ Row = 29; Doc = "73136692";   Name = "SILFREDO MARTINEZ CABARCAS";       Period = "2307"; Mora = 60000;  Salario = 1500000 },
    @{ Row = 30; Doc = "73136692";   Name = "SILFREDO MARTINEZ CABARCAS";       Period = "2308"; Mora = 60000;  Salario = 1500000 },
    @{ Row = 31; Doc = "73136692";   Name = "SILFREDO MARTINEZ CABARCAS";       Period = "2309"; Mora = 60000;  Salario = 1500000 },
    @{ Row = 32; Doc = "73136692";   Name = "SILFREDO MARTINEZ CABARCAS";       Period = "2310"; Mora = 60000;  Salario = 1500000 },
    @{ Row = 33; Doc = "73136692";   Name = "SILFREDO MARTINEZ CABARCAS";       Period = "2311"; Mora = 60000;  Salario = 1500000 },
    @{ Row = 34; Doc = "73136692";   Name = "SILFREDO MARTINEZ CABARCAS";       Period = "2312"; Mora = 60000;  Salario = 1500000 },
    @{ Row = 35; Doc = "73136692";   Name = "SILFREDO MARTINEZ CABARCAS";       Period = "2401"; Mora = 60000;  Salario = 1500000 },
    @{ Row = 36; Doc = "73136692";   Name = "SILFREDO MARTINEZ CABARCAS";       Period = "2402"; Mora = 60000;  Salario = 1500000 },
    @{ Row = 37; Doc = "73136692";   Name = "SILFREDO MARTINEZ CABARCAS";       Period = "2403"; Mora = 60000;  Salario = 1500000 },
    @{ Row = 38; Doc = "73136692";   Name = "SILFREDO MARTINEZ CABARCAS";       Period = "2404"; Mora = 46000;  Salario = 1500000 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("C$n").Value = $r.Doc
    $ws.Range("D$n").Value = $r.Name
    $ws.Range("E$n").Value = $r.Period
    $ws.Range("F$n").Value = $r.Mora
    $ws.Range("G$n").Value = $r.Salario
}
